$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the existing column width of column BY (77) for the new column BZ (78)
$ws.Range("BZ1").EntireColumn.ColumnWidth = $ws.Range("BY1").EntireColumn.ColumnWidth()

$ws.Cells.Item(1, 78).NumberFormat = "@"
$ws.Cells.Item(1, 78).Value = "2024/11/25"
$ws.Range("BY1").Copy()
$ws.Cells.Item(1, 78).PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Cells.Item(2, 78).PasteSpecial(-4122)
$ws.Cells.Item(2, 78).Value = 103.4

$ws.Range("I3").Copy()
$ws.Cells.Item(3, 78).PasteSpecial(-4122)
$ws.Cells.Item(3, 78).Value = 130.3

$ws.Range("A4").Copy()
$ws.Cells.Item(4, 78).PasteSpecial(-4122)
$ws.Cells.Item(4, 78).Value = 140.1

$ws.Range("A5").Copy()
$ws.Cells.Item(5, 78).PasteSpecial(-4122)
$ws.Cells.Item(5, 78).Value = 177

$ws.Range("A6").Copy()
$ws.Cells.Item(6, 78).PasteSpecial(-4122)
$ws.Cells.Item(6, 78).Value = 151.8

$ws.Range("A7").Copy()
$ws.Cells.Item(7, 78).PasteSpecial(-4122)
$ws.Cells.Item(7, 78).Value = 152.8

$ws.Range("A8").Copy()
$ws.Cells.Item(8, 78).PasteSpecial(-4122)
$ws.Cells.Item(8, 78).Value = 158.8

$ws.Range("A9").Copy()
$ws.Cells.Item(9, 78).PasteSpecial(-4122)
$ws.Cells.Item(9, 78).Value = 159.3

$ws.Range("D10").Copy()
$ws.Cells.Item(10, 78).PasteSpecial(-4122)
$ws.Cells.Item(10, 78).Value = 127.3

$ws.Range("A11").Copy()
$ws.Cells.Item(11, 78).PasteSpecial(-4122)
$ws.Cells.Item(11, 78).Value = 313.9

$ws.Range("A12").Copy()
$ws.Cells.Item(12, 78).PasteSpecial(-4122)
$ws.Cells.Item(12, 78).Value = 162.4

$ws.Range("H13").Copy()
$ws.Cells.Item(13, 78).PasteSpecial(-4122)
$ws.Cells.Item(13, 78).Value = 117.1

$ws.Range("A14").Copy()
$ws.Cells.Item(14, 78).PasteSpecial(-4122)
$ws.Cells.Item(14, 78).Value = 144.3

$ws.Range("A15").Copy()
$ws.Cells.Item(15, 78).PasteSpecial(-4122)
$ws.Cells.Item(15, 78).Value = 156.1

$ws.Range("L16").Copy()
$ws.Cells.Item(16, 78).PasteSpecial(-4122)
$ws.Cells.Item(16, 78).Value = 121.6

$ws.Range("C17").Copy()
$ws.Cells.Item(17, 78).PasteSpecial(-4122)
$ws.Cells.Item(17, 78).Value = 130.1

$ws.Range("A18").Copy()
$ws.Cells.Item(18, 78).PasteSpecial(-4122)
$ws.Cells.Item(18, 78).Value = 290.9

$ws.Range("G19").Copy()
$ws.Cells.Item(19, 78).PasteSpecial(-4122)
$ws.Cells.Item(19, 78).Value = 114.3

$ws.Range("A20").Copy()
$ws.Cells.Item(20, 78).PasteSpecial(-4122)
$ws.Cells.Item(20, 78).Value = 187.7

$ws.Range("H21").Copy()
$ws.Cells.Item(21, 78).PasteSpecial(-4122)
$ws.Cells.Item(21, 78).Value = 120.1

$ws.Range("C22").Copy()
$ws.Cells.Item(22, 78).PasteSpecial(-4122)
$ws.Cells.Item(22, 78).Value = 134.3

$ws.Range("A23").Copy()
$ws.Cells.Item(23, 78).PasteSpecial(-4122)
$ws.Cells.Item(23, 78).Value = 179.8

$ws.Range("A24").Copy()
$ws.Cells.Item(24, 78).PasteSpecial(-4122)
$ws.Cells.Item(24, 78).Value = 167.4

$ws.Range("D25").Copy()
$ws.Cells.Item(25, 78).PasteSpecial(-4122)
$ws.Cells.Item(25, 78).Value = 133.2

$ws.Range("A26").Copy()
$ws.Cells.Item(26, 78).PasteSpecial(-4122)
$ws.Cells.Item(26, 78).Value = 147

$ws.Range("A27").Copy()
$ws.Cells.Item(27, 78).PasteSpecial(-4122)
$ws.Cells.Item(27, 78).Value = 180.7

$ws.Range("A28").Copy()
$ws.Cells.Item(28, 78).PasteSpecial(-4122)
$ws.Cells.Item(28, 78).Value = 155.4

$ws.Range("A29").Copy()
$ws.Cells.Item(29, 78).PasteSpecial(-4122)
$ws.Cells.Item(29, 78).Value = 144.6

$ws.Range("A30").Copy()
$ws.Cells.Item(30, 78).PasteSpecial(-4122)
$ws.Cells.Item(30, 78).Value = 155.5

$ws.Range("A31").Copy()
$ws.Cells.Item(31, 78).PasteSpecial(-4122)
$ws.Cells.Item(31, 78).Value = 147.4

$ws.Range("A32").Copy()
$ws.Cells.Item(32, 78).PasteSpecial(-4122)
$ws.Cells.Item(32, 78).Value = 150.6

$ws.Range("A33").Copy()
$ws.Cells.Item(33, 78).PasteSpecial(-4122)
$ws.Cells.Item(33, 78).Value = 155.1

$ws.Range("A34").Copy()
$ws.Cells.Item(34, 78).PasteSpecial(-4122)
$ws.Cells.Item(34, 78).Value = 146.6

$ws.Range("A35").Copy()
$ws.Cells.Item(35, 78).PasteSpecial(-4122)
$ws.Cells.Item(35, 78).Value = 141.9

$ws.Range("A36").Copy()
$ws.Cells.Item(36, 78).PasteSpecial(-4122)
$ws.Cells.Item(36, 78).Value = 141.6

$ws.Range("E37").Copy()
$ws.Cells.Item(37, 78).PasteSpecial(-4122)
$ws.Cells.Item(37, 78).Value = 119.4

$ws.Range("A38").Copy()
$ws.Cells.Item(38, 78).PasteSpecial(-4122)
$ws.Cells.Item(38, 78).Value = 193.4

$ws.Range("D39").Copy()
$ws.Cells.Item(39, 78).PasteSpecial(-4122)
$ws.Cells.Item(39, 78).Value = 132.5

$ws.Range("B40").Copy()
$ws.Cells.Item(40, 78).PasteSpecial(-4122)
$ws.Cells.Item(40, 78).Value = 132.4

$ws.Range("C41").Copy()
$ws.Cells.Item(41, 78).PasteSpecial(-4122)
$ws.Cells.Item(41, 78).Value = 135.5

$ws.Range("K42").Copy()
$ws.Cells.Item(42, 78).PasteSpecial(-4122)
$ws.Cells.Item(42, 78).Value = 134.7

$ws.Range("A43").Copy()
$ws.Cells.Item(43, 78).PasteSpecial(-4122)
$ws.Cells.Item(43, 78).Value = 220.8

$ws.Range("H44").Copy()
$ws.Cells.Item(44, 78).PasteSpecial(-4122)
$ws.Cells.Item(44, 78).Value = 118.7

$ws.Range("A45").Copy()
$ws.Cells.Item(45, 78).PasteSpecial(-4122)
$ws.Cells.Item(45, 78).Value = 232.6

$ws.Range("A46").Copy()
$ws.Cells.Item(46, 78).PasteSpecial(-4122)
$ws.Cells.Item(46, 78).Value = 175.4

$ws.Range("A47").Copy()
$ws.Cells.Item(47, 78).PasteSpecial(-4122)
$ws.Cells.Item(47, 78).Value = 156

$ws.Range("J48").Copy()
$ws.Cells.Item(48, 78).PasteSpecial(-4122)
$ws.Cells.Item(48, 78).Value = 133.5

$ws.Range("O49").Copy()
$ws.Cells.Item(49, 78).PasteSpecial(-4122)
$ws.Cells.Item(49, 78).Value = 114.7

$ws.Range("A50").Copy()
$ws.Cells.Item(50, 78).PasteSpecial(-4122)
$ws.Cells.Item(50, 78).Value = 163.7

$ws.Range("B51").Copy()
$ws.Cells.Item(51, 78).PasteSpecial(-4122)
$ws.Cells.Item(51, 78).Value = 136.7

$ws.Range("A52").Copy()
$ws.Cells.Item(52, 78).PasteSpecial(-4122)
$ws.Cells.Item(52, 78).Value = 174.4

$ws.Range("AB53").Copy()
$ws.Cells.Item(53, 78).PasteSpecial(-4122)
$ws.Cells.Item(53, 78).Value = 95.8

$excel.CutCopyMode = 0